$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "69.681.09"
$ws.Range("E2").Value = "  +3.10%  "

$ws.Range("D3").Value = "3.387.70"
$ws.Range("E3").Value = "  +4.24%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "190.96"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "593.93"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.37%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  +0.62%  "

$ws.Range("E9").Value = "  +2.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.77"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.16%  "

$ws.Range("E11").Value = "  +1.80%  "

$ws.Range("D12").Value = "3.978.62"
$ws.Range("E12").Value = "  +4.42%  "

$ws.Range("E13").Value = "  -0.75%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.73"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.79%  "

$ws.Range("D15").Value = "69.667.04"
$ws.Range("E15").Value = "  +3.03%  "

$ws.Range("E16").Value = "  +1.52%  "

$ws.Range("D17").Value = "3.379.52"
$ws.Range("E17").Value = "  +5.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "452.33"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +14.72%  "

$ws.Range("E19").Value = "  +1.25%  "

$ws.Range("E20").Value = "  +1.87%  "

$ws.Range("E21").Value = "  +2.99%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "76.23"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +6.72%  "

$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("E25").Value = "  +3.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.191"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.53"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.68%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("E29").Value = "  +3.36%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "23.50"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.65%  "

$ws.Range("E31").Value = "  +1.65%  "

$ws.Range("E32").Value = "  +2.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.00"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.16%  "

$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("E35").Value = "  +7.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.41"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.18%  "

$ws.Range("E37").Value = "  +2.87%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "28.35"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +6.29%  "

$ws.Range("E39").Value = "  +1.02%  "

$ws.Range("E40").Value = "  +1.52%  "

$ws.Range("E41").Value = "  +2.06%  "

$ws.Range("D42").Value = "2.754.11"
$ws.Range("E42").Value = "  +5.34%  "

$ws.Range("E43").Value = "  +1.54%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.59"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.43%  "

$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0689"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.17"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.37%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "340.08"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.80%  "

$ws.Range("E48").Value = "  +2.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.69"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +6.14%  "

$ws.Range("E50").Value = "  +5.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.34"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.06%  "
